# Auto-generated edit script: updates market-price-derived columns (H-N)
# on several rows across all 8 item-category sheets, per the scheduled
# market data refresh described in the commit.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Cells.Item(17, 8).Value = 689.7234
$ws.Cells.Item(17, 10).Value = 662.95557
$ws.Cells.Item(17, 12).Value = 1988.86671
$ws.Cells.Item(17, 14).Value = -2324.86671
# row 31
$ws.Cells.Item(31, 8).Value = 10280
$ws.Cells.Item(31, 9).Value = 373.33334
$ws.Cells.Item(31, 11).Value = 1120.00002
$ws.Cells.Item(31, 13).Value = -890.0000199999999
# row 33
$ws.Cells.Item(33, 8).Value = 236.3158
$ws.Cells.Item(33, 9).Value = 238.75
$ws.Cells.Item(33, 10).Value = 223.33333
$ws.Cells.Item(33, 11).Value = 238.75
$ws.Cells.Item(33, 12).Value = 223.33333
$ws.Cells.Item(33, 13).Value = -9.75
$ws.Cells.Item(33, 14).Value = -681.3333299999999
# row 62
$ws.Cells.Item(62, 8).Value = 1599.8334
$ws.Cells.Item(62, 9).Value = 899.6667
$ws.Cells.Item(62, 11).Value = 899.6667
$ws.Cells.Item(62, 13).Value = -275.6667
# row 64
$ws.Cells.Item(64, 8).Value = 7514.0527
$ws.Cells.Item(64, 10).Value = 8505.77
$ws.Cells.Item(64, 12).Value = 8505.77
$ws.Cells.Item(64, 14).Value = -9001.77
# row 65
$ws.Cells.Item(65, 8).Value = 1599.8334
$ws.Cells.Item(65, 9).Value = 899.6667
$ws.Cells.Item(65, 11).Value = 4498.3335
$ws.Cells.Item(65, 13).Value = -1378.3335
# row 67
$ws.Cells.Item(67, 8).Value = 7514.0527
$ws.Cells.Item(67, 10).Value = 8505.77
$ws.Cells.Item(67, 12).Value = 8505.77
$ws.Cells.Item(67, 14).Value = -10221.77
# row 87
$ws.Cells.Item(87, 8).Value = 19833.334
$ws.Cells.Item(87, 10).Value = 19833.334
$ws.Cells.Item(87, 12).Value = 19833.334
$ws.Cells.Item(87, 14).Value = -22329.334
# row 88
$ws.Cells.Item(88, 8).Value = 4896.8125
$ws.Cells.Item(88, 9).Value = 7799.8335
$ws.Cells.Item(88, 10).Value = 3155
$ws.Cells.Item(88, 11).Value = 7799.8335
$ws.Cells.Item(88, 12).Value = 3155
$ws.Cells.Item(88, 13).Value = -7393.8335
$ws.Cells.Item(88, 14).Value = -3967
# row 90
$ws.Cells.Item(90, 8).Value = 19833.334
$ws.Cells.Item(90, 10).Value = 19833.334
$ws.Cells.Item(90, 12).Value = 59500.00199999999
$ws.Cells.Item(90, 14).Value = -71980.00199999999
# row 91
$ws.Cells.Item(91, 8).Value = 4896.8125
$ws.Cells.Item(91, 9).Value = 7799.8335
$ws.Cells.Item(91, 10).Value = 3155
$ws.Cells.Item(91, 11).Value = 7799.8335
$ws.Cells.Item(91, 12).Value = 3155
$ws.Cells.Item(91, 13).Value = -6395.8335
$ws.Cells.Item(91, 14).Value = -5963
# row 125
$ws.Cells.Item(125, 8).Value = 836.0909
$ws.Cells.Item(125, 9).Value = 1244.4
$ws.Cells.Item(125, 10).Value = 495.83334
$ws.Cells.Item(125, 11).Value = 11199.6
$ws.Cells.Item(125, 12).Value = 4462.50006
$ws.Cells.Item(125, 13).Value = -8739.6
$ws.Cells.Item(125, 14).Value = -9382.50006
# row 132
$ws.Cells.Item(132, 8).Value = 20135.1
$ws.Cells.Item(132, 9).Value = 2064.7917
$ws.Cells.Item(132, 11).Value = 6194.375100000001
$ws.Cells.Item(132, 13).Value = -3664.375100000001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Cells.Item(45, 8).Value = 6323.6816
$ws.Cells.Item(45, 9).Value = 4069.2222
$ws.Cells.Item(45, 10).Value = 7884.4614
$ws.Cells.Item(45, 11).Value = 4069.2222
$ws.Cells.Item(45, 12).Value = 7884.4614
$ws.Cells.Item(45, 13).Value = -3692.2222
$ws.Cells.Item(45, 14).Value = -8638.4614
# row 61
$ws.Cells.Item(61, 8).Value = 5201.4
$ws.Cells.Item(61, 9).Value = 2140
$ws.Cells.Item(61, 10).Value = 5467.609
$ws.Cells.Item(61, 11).Value = 2140
$ws.Cells.Item(61, 12).Value = 5467.609
$ws.Cells.Item(61, 13).Value = -1928
$ws.Cells.Item(61, 14).Value = -5891.609
# row 63
$ws.Cells.Item(63, 8).Value = 6040.852
$ws.Cells.Item(63, 10).Value = 8999.700000000001
$ws.Cells.Item(63, 12).Value = 8999.700000000001
$ws.Cells.Item(63, 14).Value = -10371.7
# row 66
$ws.Cells.Item(66, 8).Value = 6040.852
$ws.Cells.Item(66, 10).Value = 8999.700000000001
$ws.Cells.Item(66, 12).Value = 44998.5
$ws.Cells.Item(66, 14).Value = -51862.5
# row 102
$ws.Cells.Item(102, 8).Value = 10755673
$ws.Cells.Item(102, 9).Value = 3080.96
$ws.Cells.Item(102, 11).Value = 3080.96
$ws.Cells.Item(102, 13).Value = -1458.96
# row 122
$ws.Cells.Item(122, 8).Value = 4447.305
$ws.Cells.Item(122, 9).Value = 3765.3438
$ws.Cells.Item(122, 11).Value = 11296.0314
$ws.Cells.Item(122, 13).Value = -8846.0314
# row 123
$ws.Cells.Item(123, 8).Value = 30357.143
$ws.Cells.Item(123, 10).Value = 30357.143
$ws.Cells.Item(123, 12).Value = 30357.143
$ws.Cells.Item(123, 14).Value = -40157.143
# row 132
$ws.Cells.Item(132, 8).Value = 1716.6842
$ws.Cells.Item(132, 9).Value = 1766.0588
$ws.Cells.Item(132, 10).Value = 1297
$ws.Cells.Item(132, 11).Value = 5298.1764
$ws.Cells.Item(132, 12).Value = 3891
$ws.Cells.Item(132, 13).Value = -2768.1764
$ws.Cells.Item(132, 14).Value = -8951
# row 136
$ws.Cells.Item(136, 8).Value = 5201.4
$ws.Cells.Item(136, 9).Value = 2140
$ws.Cells.Item(136, 10).Value = 5467.609
$ws.Cells.Item(136, 11).Value = 6420
$ws.Cells.Item(136, 12).Value = 16402.827
$ws.Cells.Item(136, 13).Value = -3870
$ws.Cells.Item(136, 14).Value = -21502.827

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Cells.Item(86, 8).Value = 30945.441
$ws.Cells.Item(86, 9).Value = 37141.145
$ws.Cells.Item(86, 10).Value = 2032.1666
$ws.Cells.Item(86, 11).Value = 37141.145
$ws.Cells.Item(86, 12).Value = 2032.1666
$ws.Cells.Item(86, 13).Value = -36018.145
$ws.Cells.Item(86, 14).Value = -4278.1666
# row 89
$ws.Cells.Item(89, 8).Value = 30945.441
$ws.Cells.Item(89, 9).Value = 37141.145
$ws.Cells.Item(89, 10).Value = 2032.1666
$ws.Cells.Item(89, 11).Value = 185705.725
$ws.Cells.Item(89, 12).Value = 10160.833
$ws.Cells.Item(89, 13).Value = -180089.725
$ws.Cells.Item(89, 14).Value = -21392.833
# row 107
$ws.Cells.Item(107, 8).Value = 12602.111
$ws.Cells.Item(107, 10).Value = 7182.6
$ws.Cells.Item(107, 12).Value = 7182.6
$ws.Cells.Item(107, 14).Value = -11022.6
# row 134
$ws.Cells.Item(134, 8).Value = 3106.0625
$ws.Cells.Item(134, 9).Value = 3106.0625
$ws.Cells.Item(134, 11).Value = 9318.1875
$ws.Cells.Item(134, 13).Value = -6783.1875

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Cells.Item(7, 8).Value = 261.625
$ws.Cells.Item(7, 9).Value = 162.75
$ws.Cells.Item(7, 10).Value = 311.0625
$ws.Cells.Item(7, 11).Value = 162.75
$ws.Cells.Item(7, 12).Value = 311.0625
$ws.Cells.Item(7, 13).Value = -49.75
$ws.Cells.Item(7, 14).Value = -537.0625
# row 31
$ws.Cells.Item(31, 8).Value = 2680.087
$ws.Cells.Item(31, 9).Value = 2332.35
$ws.Cells.Item(31, 11).Value = 2332.35
$ws.Cells.Item(31, 13).Value = -2037.35
# row 34
$ws.Cells.Item(34, 8).Value = 2680.087
$ws.Cells.Item(34, 9).Value = 2332.35
$ws.Cells.Item(34, 11).Value = 2332.35
$ws.Cells.Item(34, 13).Value = -2130.35
# row 99
$ws.Cells.Item(99, 8).Value = 2193.375
$ws.Cells.Item(99, 10).Value = 2399
$ws.Cells.Item(99, 12).Value = 2399
$ws.Cells.Item(99, 14).Value = -5395
# row 126
$ws.Cells.Item(126, 8).Value = 2193.375
$ws.Cells.Item(126, 10).Value = 2399
$ws.Cells.Item(126, 12).Value = 7197
$ws.Cells.Item(126, 14).Value = -12137
# row 132
$ws.Cells.Item(132, 8).Value = 3712.4443
$ws.Cells.Item(132, 9).Value = 3880.2856
$ws.Cells.Item(132, 10).Value = 3125
$ws.Cells.Item(132, 11).Value = 11640.8568
$ws.Cells.Item(132, 12).Value = 9375
$ws.Cells.Item(132, 13).Value = -9110.856800000001
$ws.Cells.Item(132, 14).Value = -14435
# row 134
$ws.Cells.Item(134, 8).Value = 4228.973
$ws.Cells.Item(134, 9).Value = 4145.853
$ws.Cells.Item(134, 11).Value = 12437.559
$ws.Cells.Item(134, 13).Value = -9902.559000000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 113
$ws.Cells.Item(113, 8).Value = 1051.125
$ws.Cells.Item(113, 10).Value = 768
$ws.Cells.Item(113, 12).Value = 2304
$ws.Cells.Item(113, 14).Value = -6644

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 20
$ws.Cells.Item(20, 8).Value = 7067230
$ws.Cells.Item(20, 10).Value = 34459.8
$ws.Cells.Item(20, 12).Value = 34459.8
$ws.Cells.Item(20, 14).Value = -34949.8
# row 35
$ws.Cells.Item(35, 8).Value = 29369.572
$ws.Cells.Item(35, 9).Value = 28513.6
$ws.Cells.Item(35, 10).Value = 31509.5
$ws.Cells.Item(35, 11).Value = 28513.6
$ws.Cells.Item(35, 12).Value = 31509.5
$ws.Cells.Item(35, 13).Value = -28215.6
$ws.Cells.Item(35, 14).Value = -32105.5
# row 80
$ws.Cells.Item(80, 8).Value = 27873936
$ws.Cells.Item(80, 9).Value = 142985.5
$ws.Cells.Item(80, 11).Value = 142985.5
$ws.Cells.Item(80, 13).Value = -141987.5
# row 83
$ws.Cells.Item(83, 8).Value = 27873936
$ws.Cells.Item(83, 9).Value = 142985.5
$ws.Cells.Item(83, 11).Value = 714927.5
$ws.Cells.Item(83, 13).Value = -709935.5
# row 102
$ws.Cells.Item(102, 8).Value = 2912.7407
$ws.Cells.Item(102, 9).Value = 2365.76
$ws.Cells.Item(102, 11).Value = 2365.76
$ws.Cells.Item(102, 13).Value = -743.7600000000002
# row 113
$ws.Cells.Item(113, 8).Value = 6354.346
$ws.Cells.Item(113, 9).Value = 4185.05
$ws.Cells.Item(113, 10).Value = 13585.333
$ws.Cells.Item(113, 11).Value = 4185.05
$ws.Cells.Item(113, 12).Value = 13585.333
$ws.Cells.Item(113, 13).Value = -2015.05
$ws.Cells.Item(113, 14).Value = -17925.333
# row 122
$ws.Cells.Item(122, 8).Value = 4063.5
$ws.Cells.Item(122, 9).Value = 3343.5625
$ws.Cells.Item(122, 11).Value = 10030.6875
$ws.Cells.Item(122, 13).Value = -7580.6875

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 59
$ws.Cells.Item(59, 8).Value = 36750.75
$ws.Cells.Item(59, 10).Value = 36750.75
$ws.Cells.Item(59, 12).Value = 36750.75
$ws.Cells.Item(59, 14).Value = -38058.75
# row 95
$ws.Cells.Item(95, 8).Value = 60000
$ws.Cells.Item(95, 10).Value = 60000
$ws.Cells.Item(95, 12).Value = 60000
$ws.Cells.Item(95, 14).Value = -65492
# row 132
$ws.Cells.Item(132, 8).Value = 5164.6733
$ws.Cells.Item(132, 9).Value = 4493.3823
$ws.Cells.Item(132, 11).Value = 13480.1469
$ws.Cells.Item(132, 13).Value = -10950.1469

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 47
$ws.Cells.Item(47, 8).Value = 34179
$ws.Cells.Item(47, 9).Value = 34179
$ws.Cells.Item(47, 11).Value = 34179
$ws.Cells.Item(47, 13).Value = -33607
# row 81
$ws.Cells.Item(81, 8).Value = 15105010
$ws.Cells.Item(81, 9).Value = 16668100
$ws.Cells.Item(81, 10).Value = 13932693
$ws.Cells.Item(81, 11).Value = 33336200
$ws.Cells.Item(81, 12).Value = 27865386
$ws.Cells.Item(81, 13).Value = -33335139
$ws.Cells.Item(81, 14).Value = -27867508
# row 84
$ws.Cells.Item(84, 8).Value = 15105010
$ws.Cells.Item(84, 9).Value = 16668100
$ws.Cells.Item(84, 10).Value = 13932693
$ws.Cells.Item(84, 11).Value = 166681000
$ws.Cells.Item(84, 12).Value = 139326930
$ws.Cells.Item(84, 13).Value = -166675696
$ws.Cells.Item(84, 14).Value = -139337538
# row 95
$ws.Cells.Item(95, 8).Value = 34210
$ws.Cells.Item(95, 10).Value = 34210
$ws.Cells.Item(95, 12).Value = 34210
$ws.Cells.Item(95, 14).Value = -39702
# row 122
$ws.Cells.Item(122, 8).Value = 1300.25
$ws.Cells.Item(122, 9).Value = 1192
$ws.Cells.Item(122, 11).Value = 3576
$ws.Cells.Item(122, 13).Value = -1126
# row 136
$ws.Cells.Item(136, 8).Value = 3941.6667
$ws.Cells.Item(136, 9).Value = 3987.8647
$ws.Cells.Item(136, 10).Value = 3599.8
$ws.Cells.Item(136, 11).Value = 11963.5941
$ws.Cells.Item(136, 12).Value = 10799.4
$ws.Cells.Item(136, 13).Value = -9413.5941
$ws.Cells.Item(136, 14).Value = -15899.4

